# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Values that look numeric (e.g. "1.002", "23.136.14") are prefixed with a
# leading apostrophe so Excel stores them as text, matching the source
# sheet's inline-string cells instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.136.14"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "'1.602.44"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'301.38"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "'0.3769"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("D8").Value = "'0.3655"
$ws.Range("E8").Value = "  -4.53%  "
$ws.Range("D9").Value = "'48.01"
$ws.Range("E9").Value = "  -6.37%  "
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "'1.275"
$ws.Range("E11").Value = "  -5.88%  "
$ws.Range("D12").Value = "'0.08064"
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("D13").Value = "'22.95"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("D14").Value = "'6.619"
$ws.Range("E14").Value = "  -7.04%  "
$ws.Range("D15").Value = "'7.657"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").Value = "'0.00001265"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "'1.594.88"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").Value = "'91.46"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "'0.06797"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'18.37"
$ws.Range("E20").Value = "  -7.07%  "
$ws.Range("D21").Value = "'6.586"
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'12.99"
$ws.Range("E23").Value = "  -4.60%  "
$ws.Range("D24").Value = "'23.146.07"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").Value = "'2.357"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("D26").Value = "'2.899"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("D27").Value = "'21.05"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").Value = "'150.59"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "'5.264"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("D30").Value = "'131.83"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("D31").Value = "'2.434"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'6.939"
$ws.Range("E32").Value = "  -10.39%  "
$ws.Range("D33").Value = "'1.771.54"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "'0.9935"
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("D35").Value = "'0.07723"
$ws.Range("E35").Value = "  -5.02%  "
$ws.Range("D36").Value = "'0.02780"
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").Value = "'6.285"
$ws.Range("E37").Value = "  -6.27%  "
$ws.Range("D38").Value = "'0.2541"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'10.10"
$ws.Range("E39").Value = "  -6.72%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.08864"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'1.392"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'0.7155"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").Value = "'12.76"
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("D44").Value = "'15.79"
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("D45").Value = "'0.6614"
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").Value = "'2.311"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").Value = "'0.07993"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("D50").Value = "'131.54"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").Value = "'1.172"
$ws.Range("E51").Value = "  -4.15%  "
